$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix column: move the "Aluno"/"Contratante" value in M5 and M6 over to column N
$ws.Range("N5").Value = $ws.Range("M5").Value2
$ws.Range("M5").ClearContents()

$ws.Range("N6").Value = $ws.Range("M6").Value2
$ws.Range("M6").ClearContents()

# 2. Remove the erroneous/duplicate data row (old row 8) entirely
$ws.Range("A8:R8").ClearContents()

# 3. Always place formatted (but empty) date cells 10 rows below the last
#    data row, through row 32, ready to receive the NPS report data.
$ws.Range("C7").Copy()
$ws.Range("C8:D32").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 4. Update the selection / view state
$ws.Range("B10").Select()
